# TS2500 now has checks and a write_action operation.
# Add the 8 new to-do items to the priorities sheet, then re-sort the
# table (A2:E74) descending by the computed Priority column (E), matching
# the workbook's existing autofilter/sort-state behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name, Importance (B), Difficulty (C), Time in minutes (D)
# Priority (E) = B*C/D, entered as a formula like the existing rows.
# NB: the order below matters - it controls the order new entries are
# appended to the shared-string table.
$newItems = @(
    @("mK and uK self-heating checks", 1, 3, 1),
    @("Uncertainty in instrument files", 1, 4, 1),
    @("README.md", 1, 3, 1),
    @("Validate software", 1, 3, 3),
    @("2900 fan speed?", 1, 1, 3),
    @("Remind if instrument calibration due", 1, 4, 2),
    @("Version number", 1, 2, 1),
    @("Document tracability", 1, 5, 1)
)

$lastRow = 66
$row = $lastRow + 1
foreach ($item in $newItems) {
    $ws.Cells.Item($row, 1).Value2 = $item[0]
    $ws.Cells.Item($row, 2).Value2 = $item[1]
    $ws.Cells.Item($row, 3).Value2 = $item[2]
    $ws.Cells.Item($row, 4).Value2 = $item[3]
    $ws.Cells.Item($row, 5).Formula = "=B$row*C$row/D$row"
    $row = $row + 1
}

$newLastRow = $row - 1

# Re-sort the full data body (still below the header row) descending on
# the Priority column, same as the sheet's existing autoFilter sort.
$sortRange = $ws.Range("A2:E$newLastRow")
$keyRange = $ws.Range("E2:E$newLastRow")
$sortRange.Sort($keyRange, 2)

# Match the cursor landing on E1 after the resort/refilter.
$ws.Range("E1").Select()
